$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for the two new columns, styled like the other headers (B1:H1)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows 2-22: I column (9) is always 1, J column (10) mirrors column H (8)
for ($row = 2; $row -le 22; $row++) {
    $ws.Cells.Item($row, 9).Value = 1
    $ws.Cells.Item($row, 10).Value = $ws.Cells.Item($row, 8).Value2
}

# Row 23 breaks the pattern - explicit values
$ws.Cells.Item(23, 9).Value = 6
$ws.Cells.Item(23, 10).Value = 7
